$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, exactly as in the source data,
# so values like "0.0950" or "20.70" keep trailing zeros instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.533.57'
$ws.Range("D3").Value = '2.240.66'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '245.14'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Value = '0.629'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '75.08'
$ws.Range("E7").Value = '  -2.53%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").Value = '43.73'
$ws.Range("E10").Value = '  +4.77%  '
$ws.Range("D11").Value = '0.0950'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = '14.52'
$ws.Range("E14").Value = '  -2.60%  '
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '2.233.25'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '42.280.69'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '0.0000104'
$ws.Range("E18").Value = '  +6.05%  '
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '71.99'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '10.63'
$ws.Range("E21").Value = '  +47.04%  '
$ws.Range("D22").Value = '231.23'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("E23").Value = '  -5.30%  '
$ws.Range("D24").Value = '11.71'
$ws.Range("E24").Value = '  +2.51%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  +4.14%  '
$ws.Range("D29").Value = '166.86'
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = '20.70'
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '5.98'
$ws.Range("E31").Value = '  +22.02%  '
$ws.Range("D32").Value = '0.0817'
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = '30.03'
$ws.Range("E34").Value = '  -11.06%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").Value = '4.61'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").Value = '0.0312'
$ws.Range("E37").Value = '  +3.36%  '
$ws.Range("E38").Value = '  -6.11%  '
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("D40").Value = '5.71'
$ws.Range("E40").Value = '  -4.11%  '
$ws.Range("D41").Value = '63.58'
$ws.Range("E41").Value = '  +3.73%  '
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("D43").Value = '106.56'
$ws.Range("E43").Value = '  -6.09%  '
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("D45").Value = '0.103'
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("E47").Value = '  +5.53%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("E51").Value = '  +1.36%  '
